$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$table = $ws.ListObjects.Item("Tabelle1")

# Insert a new blank row right after the header (row 1), which shifts the
# existing data rows (2-18) down to (3-19).
$ws.Rows(2).Insert()

# Fill in the new row with the login-endpoint entry.
$ws.Range("A2").Value = "POST"
$ws.Range("B2").Value = "/user/login"
$ws.Range("C2").Value = "Loginfunktion, welche das Passwortwort überprüft"
$ws.Range("D2").Value = "fertig"
$ws.Range("E2").Value = "getestet"

# The table's own range doesn't auto-grow on a plain row insert, so extend
# it explicitly to include the new row.
$table.Resize($ws.Range("A1:E19"))

# Move the selection to E2, matching the author's cursor position after editing.
$ws.Range("E2").Select()
